# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Update column G (K) values for rows 2-24 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 7
    3  = 2
    4  = 2
    5  = 6
    6  = 5
    7  = 1
    8  = 5
    9  = 6
    10 = 2
    11 = 3
    12 = 3
    13 = 4
    14 = 4
    15 = 3
    16 = 0
    17 = 2
    18 = 1
    19 = 5
    20 = 1
    21 = 2
    22 = 2
    23 = 2
    24 = 1
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
